$d = $word.ActiveDocument

# The last paragraph currently reads:
#   "상속을 통해 다른 유형의 좀비, 솔져 추가"
# split across runs:
#   R1 (eastAsia hint): "상속을 통해 다른 유형의 좀비,"
#   R2 (no rPr):         " "
#   proofErr spellStart
#   R3 (eastAsia hint):  "솔져"
#   proofErr spellEnd
#   R4 (eastAsia hint):  " 추가"
#
# Target text is "상속을 통한 Object 관리" split as:
#   T1 (eastAsia hint): "상속을 통한 "
#   T2 (no rPr):         "Object "
#   T3 (eastAsia hint):  "관리"
#
# Plus: a brand-new empty list paragraph is appended right after it.

$lastPara = $d.Paragraphs.Last
$lastStart = $lastPara.Range.Start

# --- Step 1: collapse " " + [spellStart] + "솔져" + [spellEnd] + the
# leading space of " 추가" into a single plain run "Object ". Swallowing
# that extra leading space (consumed from R4) is what lets the engine
# drop both now-orphaned proofErr markers; R4 keeps its remaining text
# ("추가") and its original eastAsia-hinted run.
$objRange = $d.Range($lastStart + 17, $lastStart + 21)
$objRange.Text = "Object "

# --- Step 2: turn the remaining "추가" run into "관리" (still inside the
# original eastAsia-hinted run, so the hint is preserved).
$tailRange = $lastPara.Range.Duplicate
$tailRange.Find.Execute("추가", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailRange.Text = "관리"

# --- Step 3: fix up the first run's text (still the eastAsia-hinted run,
# so the hint is preserved); keep the trailing space that now separates
# it from "Object ".
$d.Content.Find.Execute("상속을 통해 다른 유형의 좀비,", $true, $false, $false, $false, $false, $true, 1, $false, "상속을 통한 ", 2) | Out-Null

# --- Step 4: append a brand-new, completely empty paragraph right after
# this one, sharing the same list paragraph formatting (pStyle a3 /
# numPr ilvl0,numId1 / ind leftChars0) but with no runs inside it.
$lastPara2 = $d.Paragraphs.Last
$para2Range = $lastPara2.Range
$para2Len = $para2Range.Text.Length
$insertPos = $para2Range.Start + $para2Len - 1
$insertionPoint = $d.Range($insertPos, $insertPos)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:leftChars="0"/></w:pPr></w:p>'
$insertionPoint.InsertXML($newParaXml)
